$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A27:E27").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)
Write-Output "done"
